# Hortaliza, Agrícola del Norte S.A. de Arica - Cebollín baby
# Weekly price update: two new daily price observations are inserted into
# the historical series (one near the top of the range, one a few rows
# later), shifting all subsequent rows down. No existing data values change
# - they simply move to new row numbers as the new rows are inserted above
# them.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert new row at position 17 -----------------------------------
$ws.Rows(17).Insert()

$ws.Range("A17").Value = 1
$ws.Range("B17").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C17").Value = "Arica y Parinacota"
$ws.Range("D17").Value = 44970
$ws.Range("E17").Value = 15
$ws.Range("F17").Value = 100112038
$ws.Range("G17").Value = "Cebollín baby"
$ws.Range("H17").Value = "Sin especificar"
$ws.Range("I17").Value = "Primera"
$ws.Range("J17").Value = 200
$ws.Range("K17").Value = 4500
$ws.Range("L17").Value = 5000
$ws.Range("M17").Value = 4750
$ws.Range("N17").Value = "$/paquete 1,5 a 2 kilos"
$ws.Range("O17").Value = "Región de Arica y Parinacota"
$ws.Range("P17").Value = 2375
$ws.Range("Q17").Value = 2
$ws.Range("R17").Value = "Hortaliza"

# --- Insert second new row at position 22 -----------------------------
# (rows 18-21 now hold what were rows 17-20; row 21 now holds what was
# row 20, so inserting here pushes the old row 21 onward down by one more)
$ws.Rows(22).Insert()

$ws.Range("A22").Value = 1
$ws.Range("B22").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C22").Value = "Arica y Parinacota"
$ws.Range("D22").Value = 44971
$ws.Range("E22").Value = 15
$ws.Range("F22").Value = 100112038
$ws.Range("G22").Value = "Cebollín baby"
$ws.Range("H22").Value = "Sin especificar"
$ws.Range("I22").Value = "Segunda"
$ws.Range("J22").Value = 300
$ws.Range("K22").Value = 2500
$ws.Range("L22").Value = 2800
$ws.Range("M22").Value = 2600
$ws.Range("N22").Value = "$/paquete 1,5 a 2 kilos"
$ws.Range("O22").Value = "Región de Arica y Parinacota"
$ws.Range("P22").Value = 1300
$ws.Range("Q22").Value = 2
$ws.Range("R22").Value = "Hortaliza"
